$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '35.534.10'
$ws.Range("E2").Value = '  +1.33%  '

# Row 3
$ws.Range("D3").Value = '1.907.93'
$ws.Range("E3").Value = '  +2.68%  '

# Row 4
$ws.Range("E4").Value = '  +0.39%  '

# Row 5
$ws.Range("D5").Value = "'246.53"
$ws.Range("E5").Value = '  +3.62%  '

# Row 6
$ws.Range("D6").Value = "'0.664"
$ws.Range("E6").Value = '  +6.67%  '

# Row 7
$ws.Range("E7").Value = '  +0.30%  '

# Row 8
$ws.Range("D8").Value = "'41.96"
$ws.Range("E8").Value = '  -1.58%  '

# Row 9
$ws.Range("D9").Value = "'0.344"
$ws.Range("E9").Value = '  +4.15%  '

# Row 10
$ws.Range("D10").Value = "'49.15"
$ws.Range("E10").Value = '  +5.47%  '

# Row 11
$ws.Range("E11").Value = '  +2.52%  '

# Row 12
$ws.Range("D12").Value = "'0.100"
$ws.Range("E12").Value = '  +0.99%  '

# Row 13
$ws.Range("D13").Value = '2.184.01'
$ws.Range("E13").Value = '  +2.65%  '

# Row 14
$ws.Range("D14").Value = "'12.36"
$ws.Range("E14").Value = '  +8.18%  '

# Row 15
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '1.925.19'
$ws.Range("E15").Value = '  +3.52%  '

# Row 16
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").Value = "'0.700"
$ws.Range("E16").Value = '  +3.30%  '

# Row 17
$ws.Range("D17").Value = "'4.86"
$ws.Range("E17").Value = '  +3.20%  '

# Row 18
$ws.Range("D18").Value = '35.565.83'
$ws.Range("E18").Value = '  +1.46%  '

# Row 19
$ws.Range("D19").Value = "'72.35"
$ws.Range("E19").Value = '  +2.96%  '

# Row 20
$ws.Range("E20").Value = '  +4.63%  '

# Row 21
$ws.Range("D21").Value = "'243.92"
$ws.Range("E21").Value = '  +1.27%  '

# Row 22
$ws.Range("E22").Value = '  +3.96%  '

# Row 23
$ws.Range("D23").Value = "'4.85"
$ws.Range("E23").Value = '  +2.39%  '

# Row 24
$ws.Range("E24").Value = '  +0.32%  '

# Row 25
$ws.Range("E25").Value = '  +1.31%  '

# Row 26
$ws.Range("E26").Value = '  +14.23%  '

# Row 27
$ws.Range("D27").Value = "'171.54"
$ws.Range("E27").Value = '  +0.15%  '

# Row 28
$ws.Range("D28").Value = "'8.52"
$ws.Range("E28").Value = '  +7.50%  '

# Row 29
$ws.Range("D29").Value = "'18.28"
$ws.Range("E29").Value = '  +3.19%  '

# Row 30
$ws.Range("D30").Value = "'0.130"
$ws.Range("E30").Value = '  +4.16%  '

# Row 31
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = "'4.17"
$ws.Range("E31").Value = '  +4.32%  '

# Row 32
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = "'0.965"
$ws.Range("E32").Value = '  +22.32%  '

# Row 33
$ws.Range("E33").Value = '  +2.37%  '

# Row 34
$ws.Range("D34").Value = "'4.24"
$ws.Range("E34").Value = '  +5.53%  '

# Row 35
$ws.Range("E35").Value = '  +0.45%  '

# Row 36
$ws.Range("D36").Value = "'1.73"
$ws.Range("E36").Value = '  +7.03%  '

# Row 37
$ws.Range("E37").Value = '  -0.27%  '

# Row 38
$ws.Range("D38").Value = "'1.32"
$ws.Range("E38").Value = '  +2.34%  '

# Row 39
$ws.Range("E39").Value = '  +2.30%  '

# Row 40
$ws.Range("D40").Value = "'92.40"
$ws.Range("E40").Value = '  +0.62%  '

# Row 41
$ws.Range("E41").Value = '  +1.25%  '

# Row 42
$ws.Range("D42").Value = "'0.0635"
$ws.Range("E42").Value = '  +17.18%  '

# Row 43
$ws.Range("D43").Value = "'15.68"
$ws.Range("E43").Value = '  +5.19%  '

# Row 44
$ws.Range("D44").Value = '1.348.87'
$ws.Range("E44").Value = '  -0.48%  '

# Row 45
$ws.Range("E45").Value = '  +2.18%  '

# Row 46
$ws.Range("D46").Value = "'47.96"
$ws.Range("E46").Value = '  +40.13%  '

# Row 47
$ws.Range("B47").Value = 'HuobiToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D47").Value = "'2.42"
$ws.Range("E47").Value = '  +0.38%  '

# Row 48
$ws.Range("B48").Value = 'Gas'
$ws.Range("C48").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D48").Value = "'12.55"
$ws.Range("E48").Value = '  -3.00%  '

# Row 49
$ws.Range("D49").Value = "'2.77"
$ws.Range("E49").Value = '  +0.42%  '

# Row 50
$ws.Range("D50").Value = "'6.60"
$ws.Range("E50").Value = '  +2.53%  '

# Row 51
$ws.Range("D51").Value = '2.095.04'
$ws.Range("E51").Value = '  +2.73%  '
